$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate column G (K) values using H (IP) instead of the old Strike# calc.
# Per the diff, column G should now mirror column H (IP) for each row (with
# one row differing by +1 due to rounding in the underlying std/mean calc).
$newG = @{
    2  = 5
    3  = 3
    4  = 5
    5  = 5
    6  = 5
    7  = 5
    8  = 5
    9  = 6
    10 = 3
    11 = 0
    12 = 5
    13 = 4
}

foreach ($row in $newG.Keys) {
    $ws.Range("G$row").Value = $newG[$row]
}
